# Lecture_004.pptx - "Add files via upload" edit replay
#
# Semantic changes being applied:
#  1. The auto-updating "datetimeFigureOut" footer field cached on the
#     Slide Master and on every one of the 11 slide layouts is refreshed
#     from 10/4/2023 -> 10/12/2023 (PowerPoint re-stamps this cached text
#     whenever the deck is re-saved on a later date).
#  2. On slide 2, the table cell that used to read "hello2.py" is cleared
#     out, and the table's two rows grow to their new (taller) heights.
#  3. On slide 3, the table cell that used to read "hello3.py" is cleared
#     out (row heights on that table are unaffected).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached date field text everywhere it is stored: once on
#    the Slide Master, and once per Custom Layout.
# ---------------------------------------------------------------------
function Set-DatePlaceholderText {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = "10/12/2023"
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 2 ("hello2.py" table): clear the cell text and grow the rows.
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shape = $slide2.Shapes.Item($i)
    if ($shape.HasTable) {
        $table = $shape.Table
        $table.Rows.Item(1).Cells.Item(1).Shape.TextFrame.TextRange.Text = ""
        $table.Rows.Item(1).Height = 351190 / 12700
        $table.Rows.Item(2).Height = 1859360 / 12700
    }
}

# ---------------------------------------------------------------------
# 3) Slide 3 ("hello3.py" table): clear the cell text (heights unchanged).
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $shape = $slide3.Shapes.Item($i)
    if ($shape.HasTable) {
        $table = $shape.Table
        $table.Rows.Item(1).Cells.Item(1).Shape.TextFrame.TextRange.Text = ""
    }
}
